# Add files via upload
# - Adds two new date columns (Y = "7-jul", Z = "10-jul") with their data
# - Hides the now-unused intermediate columns C:S
# - Updates the active selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -------------------------------------------------
$ws.Range("Y1").Value = "7-jul"
$ws.Range("Z1").Value = "10-jul"

# --- New column Z data (column Y has no data in this workbook) ---------
$ws.Range("Z2").Value  = 0
$ws.Range("Z3").Value  = 12.114747032477784
$ws.Range("Z4").Value  = 19.896132867731442
$ws.Range("Z5").Value  = 25.39509033762949
$ws.Range("Z6").Value  = 0
$ws.Range("Z7").Value  = 17.249084502201711
$ws.Range("Z8").Value  = 9.6951537000179258
$ws.Range("Z9").Value  = 21.319303056650977
$ws.Range("Z10").Value = 22.581117190588198
$ws.Range("Z11").Value = 14.46941243756444
$ws.Range("Z12").Value = 0
$ws.Range("Z13").Value = 13.057731994777569
$ws.Range("Z14").Value = 0
$ws.Range("Z15").Value = 0
$ws.Range("Z16").Value = 23.050760072076518
$ws.Range("Z17").Value = 0
$ws.Range("Z18").Value = 0

# --- Hide the helper / intermediate columns (C through S) --------------
# Columns G:I and K:S had no explicit width before; collapse them to 0
# (matching the target layout) before hiding the whole C:S block.
$ws.Range("G1:I1").EntireColumn.ColumnWidth = 0
$ws.Range("K1:S1").EntireColumn.ColumnWidth = 0
$ws.Range("C:S").EntireColumn.Hidden = $true

# --- Update view: clear the frozen scroll column, move selection -------
$null = $ws.Range("AB5").Select()
